$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate cells in the exact order the original author typed them, so
# --- the shared-string table comes out in the same sequence as the source
# --- workbook (header row left-to-right, then row 2 Name-before-FSN, then
# --- the remaining rows FSN-before-Name).
$ws.Range("A1").Value = "FSN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

$ws.Range("B2").Value = "Shamanth"
$ws.Range("A2").Value = "FID1"
$ws.Range("C2").Value = 35
$ws.Range("D2").Value = 9986407821

$ws.Range("A3").Value = "FID2"
$ws.Range("B3").Value = "Naveen"
$ws.Range("C3").Value = 35
$ws.Range("D3").Value = 9743750743

$ws.Range("A4").Value = "FID3"
$ws.Range("B4").Value = "Ganaraj"
$ws.Range("C4").Value = 32
$ws.Range("D4").Value = 8086714071

$ws.Range("A5").Value = "FID4"
$ws.Range("B5").Value = "Madhura"
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = 9986429412

$ws.Range("A6").Value = "FID5"
$ws.Range("B6").Value = "Jayapadmini"
$ws.Range("C6").Value = 33
$ws.Range("D6").Value = 9535667372

$ws.Range("A7").Value = "FID6"
$ws.Range("B7").Value = "Asha"
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = 9483015748

$ws.Range("A8").Value = "FID7"
$ws.Range("B8").Value = "Shwetha"
$ws.Range("C8").Value = 30
$ws.Range("D8").Value = 9449330933

$ws.Range("A9").Value = "FID8"
$ws.Range("B9").Value = "Rithesh"
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 9449389449

$ws.Range("A10").Value = "FID9"
$ws.Range("B10").Value = "Naitik"
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 8296502878

$ws.Range("A11").Value = "FID10"
$ws.Range("B11").Value = "Vasudev"
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 8217894079

# --- Row 12 / column E are left blank in the data but still carry the
# --- table's font formatting (mirrors the blank trailing row/column in the
# --- source sheet), so touch them too before applying the font below.

# --- Apply the (non-default, color-less) font used throughout the table to
# --- the whole A1:E12 block in one shot - this also stamps blank cells
# --- (column E, row 12) with the style without giving them a value.
$ws.Range("A1:E12").Font.Name = "Calibri"

# --- Column widths (characters) so the saved sheet's <cols> matches the
# --- source file's pixel-snapped widths as closely as this engine's
# --- rounding allows.
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 18.5
$ws.Columns.Item(3).ColumnWidth = 13.5
$ws.Columns.Item(4).ColumnWidth = 13.8333333
$ws.Columns.Item(5).ColumnWidth = 18.8333333

# --- Selection cell + print orientation, matching the sheetView/pageSetup
# --- changes in the target file.
$null = $ws.Range("G15").Select()
$ws.PageSetup.Orientation = 1
